$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The only real content change is cell E8: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active selection on the sheet (matches cached view state in diff)
$ws.Range("E8").Select()
